# Updated cryptos list on Mon Nov 25 06:20:22 UTC 2024 with GitHub Actions
#
# Refreshes the per-coin Price (column D) and Volume(1h) (column E) text
# values on Sheet1, and re-sorts a few rows whose rank order changed
# (SuiNetwork/BitcoinCash swap at rows 23-24; the MantraDAO/Cosmos/
# Filecoin/Stacks/OKB rotation at rows 47-51) by rewriting their Coin
# name (B), Link (C), Price (D) and Volume(1h) (E) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every cell written below is scraped display text (prices/percentages
# formatted as strings, some with deliberate thousands-dot grouping like
# "98.325.76"), never a real number/formula. Column D in particular holds
# plain-looking decimals (e.g. "0.426", "42.46") that Excel would
# otherwise auto-coerce to a Double on assignment, so for those cells we
# force the Text number format first, assign the literal string, then
# restore the "Normal" style so no stray formatting lingers on the cell.
$updates = @(
    @{ Ref = 'D2'; Value = '98.325.76'; IsNumericCol = $true },
    @{ Ref = 'E2'; Value = '  +0.17%  '; IsNumericCol = $false },
    @{ Ref = 'D3'; Value = '3.389.64'; IsNumericCol = $true },
    @{ Ref = 'E3'; Value = '  -0.48%  '; IsNumericCol = $false },
    @{ Ref = 'E4'; Value = '  +0.00%  '; IsNumericCol = $false },
    @{ Ref = 'D5'; Value = '253.84'; IsNumericCol = $true },
    @{ Ref = 'E5'; Value = '  -0.27%  '; IsNumericCol = $false },
    @{ Ref = 'D6'; Value = '660.98'; IsNumericCol = $true },
    @{ Ref = 'E6'; Value = '  -0.04%  '; IsNumericCol = $false },
    @{ Ref = 'E7'; Value = '  +0.77%  '; IsNumericCol = $false },
    @{ Ref = 'D8'; Value = '0.426'; IsNumericCol = $true },
    @{ Ref = 'E8'; Value = '  -1.71%  '; IsNumericCol = $false },
    @{ Ref = 'E9'; Value = '  -1.06%  '; IsNumericCol = $false },
    @{ Ref = 'E10'; Value = '  -0.03%  '; IsNumericCol = $false },
    @{ Ref = 'D11'; Value = '3.388.81'; IsNumericCol = $true },
    @{ Ref = 'E11'; Value = '  -0.37%  '; IsNumericCol = $false },
    @{ Ref = 'E12'; Value = '  -2.66%  '; IsNumericCol = $false },
    @{ Ref = 'D13'; Value = '42.46'; IsNumericCol = $true },
    @{ Ref = 'E13'; Value = '  +0.49%  '; IsNumericCol = $false },
    @{ Ref = 'D14'; Value = '97.903.58'; IsNumericCol = $true },
    @{ Ref = 'E14'; Value = '  +0.02%  '; IsNumericCol = $false },
    @{ Ref = 'D15'; Value = '6.14'; IsNumericCol = $true },
    @{ Ref = 'E15'; Value = '  -6.11%  '; IsNumericCol = $false },
    @{ Ref = 'E16'; Value = '  -3.60%  '; IsNumericCol = $false },
    @{ Ref = 'D17'; Value = '4.014.78'; IsNumericCol = $true },
    @{ Ref = 'E17'; Value = '  -0.79%  '; IsNumericCol = $false },
    @{ Ref = 'D18'; Value = '9.22'; IsNumericCol = $true },
    @{ Ref = 'E18'; Value = '  +2.62%  '; IsNumericCol = $false },
    @{ Ref = 'D19'; Value = '3.395.79'; IsNumericCol = $true },
    @{ Ref = 'E19'; Value = '  -0.37%  '; IsNumericCol = $false },
    @{ Ref = 'E20'; Value = '  +3.46%  '; IsNumericCol = $false },
    @{ Ref = 'E21'; Value = '  -6.63%  '; IsNumericCol = $false },
    @{ Ref = 'D22'; Value = '10.98'; IsNumericCol = $true },
    @{ Ref = 'E22'; Value = '  -0.58%  '; IsNumericCol = $false },
    @{ Ref = 'B23'; Value = 'BitcoinCash'; IsNumericCol = $false },
    @{ Ref = 'C23'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; IsNumericCol = $false },
    @{ Ref = 'D23'; Value = '513.23'; IsNumericCol = $true },
    @{ Ref = 'E23'; Value = '  +0.83%  '; IsNumericCol = $false },
    @{ Ref = 'B24'; Value = 'SuiNetwork'; IsNumericCol = $false },
    @{ Ref = 'C24'; Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'; IsNumericCol = $false },
    @{ Ref = 'D24'; Value = '3.45'; IsNumericCol = $true },
    @{ Ref = 'E24'; Value = '  +0.06%  '; IsNumericCol = $false },
    @{ Ref = 'D25'; Value = '7.01'; IsNumericCol = $true },
    @{ Ref = 'E25'; Value = '  +4.11%  '; IsNumericCol = $false },
    @{ Ref = 'D26'; Value = '0.0000202'; IsNumericCol = $true },
    @{ Ref = 'E26'; Value = '  -1.96%  '; IsNumericCol = $false },
    @{ Ref = 'D27'; Value = '97.03'; IsNumericCol = $true },
    @{ Ref = 'E27'; Value = '  -3.47%  '; IsNumericCol = $false },
    @{ Ref = 'D28'; Value = '12.39'; IsNumericCol = $true },
    @{ Ref = 'E28'; Value = '  -3.60%  '; IsNumericCol = $false },
    @{ Ref = 'D29'; Value = '3.570.69'; IsNumericCol = $true },
    @{ Ref = 'E29'; Value = '  -0.70%  '; IsNumericCol = $false },
    @{ Ref = 'D30'; Value = '11.61'; IsNumericCol = $true },
    @{ Ref = 'E30'; Value = '  +0.31%  '; IsNumericCol = $false },
    @{ Ref = 'E31'; Value = '  -3.36%  '; IsNumericCol = $false },
    @{ Ref = 'D32'; Value = '0.998'; IsNumericCol = $true },
    @{ Ref = 'E32'; Value = '  +0.04%  '; IsNumericCol = $false },
    @{ Ref = 'D33'; Value = '0.189'; IsNumericCol = $true },
    @{ Ref = 'E33'; Value = '  -4.35%  '; IsNumericCol = $false },
    @{ Ref = 'E34'; Value = '  +9.26%  '; IsNumericCol = $false },
    @{ Ref = 'D35'; Value = '1.00'; IsNumericCol = $true },
    @{ Ref = 'E35'; Value = '  +0.14%  '; IsNumericCol = $false },
    @{ Ref = 'D36'; Value = '0.564'; IsNumericCol = $true },
    @{ Ref = 'E36'; Value = '  -1.62%  '; IsNumericCol = $false },
    @{ Ref = 'D37'; Value = '28.98'; IsNumericCol = $true },
    @{ Ref = 'E37'; Value = '  -2.92%  '; IsNumericCol = $false },
    @{ Ref = 'E38'; Value = '  +0.71%  '; IsNumericCol = $false },
    @{ Ref = 'E39'; Value = '  +0.05%  '; IsNumericCol = $false },
    @{ Ref = 'D40'; Value = '534.65'; IsNumericCol = $true },
    @{ Ref = 'E40'; Value = '  -0.08%  '; IsNumericCol = $false },
    @{ Ref = 'E41'; Value = '  +0.32%  '; IsNumericCol = $false },
    @{ Ref = 'E42'; Value = '  -0.09%  '; IsNumericCol = $false },
    @{ Ref = 'D43'; Value = '24.43'; IsNumericCol = $true },
    @{ Ref = 'E43'; Value = '  -1.17%  '; IsNumericCol = $false },
    @{ Ref = 'D44'; Value = '0.857'; IsNumericCol = $true },
    @{ Ref = 'E44'; Value = '  -2.26%  '; IsNumericCol = $false },
    @{ Ref = 'E45'; Value = '  -1.04%  '; IsNumericCol = $false },
    @{ Ref = 'D46'; Value = '1.75'; IsNumericCol = $true },
    @{ Ref = 'E46'; Value = '  +2.06%  '; IsNumericCol = $false },
    @{ Ref = 'B47'; Value = 'Stacks'; IsNumericCol = $false },
    @{ Ref = 'C47'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'; IsNumericCol = $false },
    @{ Ref = 'D47'; Value = '2.29'; IsNumericCol = $true },
    @{ Ref = 'E47'; Value = '  +8.92%  '; IsNumericCol = $false },
    @{ Ref = 'B48'; Value = 'MantraDAO'; IsNumericCol = $false },
    @{ Ref = 'C48'; Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'; IsNumericCol = $false },
    @{ Ref = 'D48'; Value = '3.69'; IsNumericCol = $true },
    @{ Ref = 'E48'; Value = '  -1.58%  '; IsNumericCol = $false },
    @{ Ref = 'B49'; Value = 'Cosmos'; IsNumericCol = $false },
    @{ Ref = 'C49'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; IsNumericCol = $false },
    @{ Ref = 'D49'; Value = '8.69'; IsNumericCol = $true },
    @{ Ref = 'E49'; Value = '  -3.91%  '; IsNumericCol = $false },
    @{ Ref = 'B50'; Value = 'OKB'; IsNumericCol = $false },
    @{ Ref = 'C50'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; IsNumericCol = $false },
    @{ Ref = 'D50'; Value = '56.23'; IsNumericCol = $true },
    @{ Ref = 'E50'; Value = '  +5.03%  '; IsNumericCol = $false },
    @{ Ref = 'B51'; Value = 'Filecoin'; IsNumericCol = $false },
    @{ Ref = 'C51'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; IsNumericCol = $false },
    @{ Ref = 'D51'; Value = '5.60'; IsNumericCol = $true },
    @{ Ref = 'E51'; Value = '  -3.47%  '; IsNumericCol = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    if ($u.IsNumericCol) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
